$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Cells.Item(2, 6).Value = 199
$ws.Cells.Item(4, 6).Value = 1162
$ws.Cells.Item(7, 6).Value = 33
$ws.Cells.Item(10, 6).Value = 361
$ws.Cells.Item(13, 6).Value = 316
$ws.Cells.Item(14, 6).Value = 359
$ws.Cells.Item(15, 6).Value = 33
$ws.Cells.Item(16, 6).Value = 65
$ws.Cells.Item(17, 6).Value = 516
$ws.Cells.Item(18, 6).Value = 1461
$ws.Cells.Item(19, 6).Value = 5671
$ws.Cells.Item(21, 6).Value = 1585
$ws.Cells.Item(23, 6).Value = 42
$ws.Cells.Item(24, 6).Value = 26
$ws.Cells.Item(25, 6).Value = 5186
$ws.Cells.Item(26, 6).Value = 5186
$ws.Cells.Item(27, 6).Value = 126
$ws.Cells.Item(29, 6).Value = 1519
$ws.Cells.Item(30, 6).Value = 18
$ws.Cells.Item(35, 6).Value = 3805

$ws = $wb.Worksheets.Item("演出")
$ws.Cells.Item(5, 6).Value = 153
$ws.Cells.Item(8, 6).Value = 126

$ws = $wb.Worksheets.Item("本地生活")
$ws.Cells.Item(4, 6).Value = 2141

$ws = $wb.Worksheets.Item("全部类型")
$ws.Cells.Item(4, 6).Value = 2141
$ws.Cells.Item(5, 6).Value = 199
$ws.Cells.Item(7, 6).Value = 1162
$ws.Cells.Item(10, 6).Value = 33
$ws.Cells.Item(12, 6).Value = 361
$ws.Cells.Item(14, 6).Value = 316
$ws.Cells.Item(15, 6).Value = 359
$ws.Cells.Item(16, 6).Value = 33
$ws.Cells.Item(17, 6).Value = 65
$ws.Cells.Item(21, 6).Value = 516
$ws.Cells.Item(22, 6).Value = 1461
$ws.Cells.Item(23, 6).Value = 5671
$ws.Cells.Item(25, 6).Value = 1585
$ws.Cells.Item(31, 6).Value = 5186
$ws.Cells.Item(32, 6).Value = 5186
$ws.Cells.Item(33, 6).Value = 126
$ws.Cells.Item(35, 6).Value = 1519
$ws.Cells.Item(36, 6).Value = 18
$ws.Cells.Item(47, 6).Value = 3805
